$d = $word.ActiveDocument
$p = $d.Paragraphs(1)
$p.Range.ParagraphFormat.Alignment = 3  # wdAlignParagraphJustify

$r = $d.Content
$r.Collapse(0) | Out-Null

$r.InsertAfter("24.08.2020")
$r.Bold = 1
$r.LanguageID = "lv-LV"
$r.Collapse(0) | Out-Null

$r.InsertAfter(": ")
$r.LanguageID = "lv-LV"
$r.Collapse(0) | Out-Null

$r.InsertAfter("Начало дня потратил на установку необходимых программ и их настройку, снятие ограничения на установку программ для компьютера, подписывание документов, получение карточки для дверей, создание")
$r.Collapse(0) | Out-Null

$r.InsertAfter(" Git ")
$r.LanguageID = "lv-LV"
$r.Collapse(0) | Out-Null

$r.InsertAfter("репозитория итд.. Затем занимался изучением документации об ")
$r.Collapse(0) | Out-Null

$r.InsertAfter("jBPM")
$r.LanguageID = "lv-LV"
$r.Collapse(0) | Out-Null

$r.InsertAfter(". Остаток дня потратил на установку сервера и его настройку, столкнулся с несколькими ошибками ")
$r.Collapse(0) | Out-Null

$r.InsertAfter("Java, ")
$r.LanguageID = "lv-LV"
$r.Collapse(0) | Out-Null

$r.InsertAfter("связанными с оперативной памятью, решил их, немного посмотрел возможности, которые предлагает ")
$r.Collapse(0) | Out-Null

$r.InsertAfter("jBPM")
$r.LanguageID = "lv-LV"
$r.Collapse(0) | Out-Null

$r.InsertAfter(", описание создать пока что не успел.")
$r.Collapse(0) | Out-Null
